$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.678.09"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.808.06"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.17"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.57"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0712"
$ws.Range("E10").Value = "  +7.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "2.068.73"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.11"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").Value = "1.803.92"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "34.698.91"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.69"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.48"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "0.0₃0804"
$ws.Range("E20").Value = "  +8.25%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.84"
$ws.Range("E22").Value = "  +3.76%  "
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0535"
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.66"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +3.41%  "
$ws.Range("D35").Value = "1.440.38"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.644"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0193"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.14"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.966"
$ws.Range("E40").Value = "  +7.14%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.16"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("B44").Value = "Swop.fi"
$ws.Range("C44").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "329.68"
$ws.Range("E44").Value = "  +531.09%  "
$ws.Range("E45").Value = "  +6.72%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0494"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.964.03"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +8.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.21"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("E51").Value = "  +0.05%  "
